# tc_214 - add a second "unit period" pricing-type option ($/ft/period)
# and show it as selected for a couple of rows; also tidy up a leftover
# cell style and update the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 and 4 ("PRICING TYPE" column D) now use the new
# "$/ft/period" option instead of the plain "$/period" one.
$ws.Range("D2").Value = "$/ft/period"
$ws.Range("D4").Value = "$/ft/period"

# Row 3 keeps "$/period" (unchanged), but the P_PLANE column (F) values
# for rows 3 and 4 are swapped: "Wet Storage" <-> "Dry Storage".
$ws.Range("F3").Value = "Dry Storage"
$ws.Range("F4").Value = "Wet Storage"

# Drop the stray red-font style that had been applied to the CAPACITY /
# PRICING TYPE columns (C1:D15), restoring the default "Normal" style.
$ws.Range("C1:D15").Style = "Normal"

# Update the sheet's active selection.
$ws.Range("D11").Select()
